# Edit sheet Card2 by admin
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Card2")

# O1: remove trailing space from "Serviced by "
$ws.Range("O1").Value = "Serviced by"

# O2:O13 -> "nan"
for ($r = 2; $r -le 13; $r++) {
    $ws.Cells.Item($r, 15).Value = "nan"
}

# M9 -> Arabic text with trailing space
$ws.Range("M9").Value = "تم سن الفلاتس ومعايرته "
